# Update the "取得日時" (acquisition timestamp) column for all data rows
# from "2026-02-13 18:44:20" to "2026-02-13 18:56:19".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldValue = "2026-02-13 18:44:20"
$newValue = "2026-02-13 18:56:19"

for ($row = 2; $row -le 14; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
